$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "@prefix": row 13/14 content swap (image <-> unitLength) and three new
# prefix rows appended (pixels, bindata, xmlAnnotation namespaces).
# ---------------------------------------------------------------------------
$wsPrefix = $wb.Worksheets.Item("@prefix")

$wsPrefix.Range("A13").Value = "unitLength"
$wsPrefix.Range("B13").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/Unit/UnitLength#"
$wsPrefix.Range("A14").Value = "image"
$wsPrefix.Range("B14").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/xmlannotation-body-space/image/"

$wsPrefix.Range("A15").Value = "pixels"
$wsPrefix.Range("B15").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/xmlannotation-body-space/pixels/"
$wsPrefix.Range("A16").Value = "bindata"
$wsPrefix.Range("B16").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/xmlannotation-body-space/bindata/"
$wsPrefix.Range("A17").Value = "xmlAnnotation"
$wsPrefix.Range("B17").Value = "http://www.openmicroscopy.org/rdf/2016-06/ome_core/samples/xmlannotation-body-space/xmlAnnotation/"

# copy style from the last pre-existing data row onto the new rows
$wsPrefix.Range("A14:B14").Copy() | Out-Null
$wsPrefix.Range("A15:B17").PasteSpecial(-4122) | Out-Null

$wsPrefix.Columns.Item(2).ColumnWidth = 99

# ---------------------------------------------------------------------------
# Sheet "Image": fix the naming-convention values for pixels/xmlAnnotation ids
# ---------------------------------------------------------------------------
$wsImage = $wb.Worksheets.Item("Image")

$wsImage.Range("E5").Value = "pixels:pixels0:0"
$wsImage.Range("G5").Value = "xmlAnnotation:image0"

$wsImage.Columns.Item(4).ColumnWidth = 22.142857142857142
$wsImage.Columns.Item(7).ColumnWidth = 19.428571428571427

$wsImage.Range("G5").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet "Pixels": same naming-convention fix for pixels/bindata ids
# ---------------------------------------------------------------------------
$wsPixels = $wb.Worksheets.Item("Pixels")

$wsPixels.Range("B5").Value = "pixels:pixels0:0"
$wsPixels.Range("M5").Value = "bindata:bindata1"

$wsPixels.Columns.Item(2).ColumnWidth = 13.714285714285714
$wsPixels.Columns.Item(13).ColumnWidth = 14.714285714285714

# ---------------------------------------------------------------------------
# Sheet "Binary_Data": naming-convention fix for bindata id
# ---------------------------------------------------------------------------
$wsBinary = $wb.Worksheets.Item("Binary_Data")

$wsBinary.Range("B5").Value = "bindata:bindata1"

$wsBinary.Columns.Item(2).ColumnWidth = 14.714285714285714
$wsBinary.Columns.Item(3).ColumnWidth = 10.428571428571429

# ---------------------------------------------------------------------------
# Sheet "Structured_Annotations": naming-convention fix for xmlAnnotation id
# (also fixes the "xmlAnnoation" typo)
# ---------------------------------------------------------------------------
$wsStruct = $wb.Worksheets.Item("Structured_Annotations")

$wsStruct.Range("C5").Value = "xmlAnnotation:image0"

$wsStruct.Columns.Item(3).ColumnWidth = 20.142857142857142

# ---------------------------------------------------------------------------
# Sheet "XML_Annotation": naming-convention fix for xmlAnnotation id
# ---------------------------------------------------------------------------
$wsXml = $wb.Worksheets.Item("XML_Annotation")

$wsXml.Range("B5").Value = "xmlAnnotation:image0"

# ---------------------------------------------------------------------------
# Restore "@prefix" as the active/selected sheet (it was active before the
# edits and must remain so - only the Image sheet's own selection moves).
# ---------------------------------------------------------------------------
$wsPrefix.Activate() | Out-Null
